# Commit: "Create the data model, create the database context"
#
# The underlying edit adds a brand-new worksheet ("Лист1") to the workbook,
# cloned from the existing "Nädal 2" sheet (same column widths, row
# formatting, cell styles and merged-cell layout), then:
#   - updates the header date (G4)
#   - keeps only the first data row (row 7) populated, with new values
#   - clears out the remaining data rows (8-16)
#   - removes the leftover cell comment that "Nädal 2" carried on H14
#     (it doesn't apply any more since that cell's data was wiped)
#   - makes the new sheet the active / selected tab
#   - the previously-active "Nädal 2" tab gets its whole-sheet selected
#     (Ctrl+A style selection) and loses the "last active" tab flag

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Nädal 2")

# Clone "Nädal 2" (preserves column widths / styles / merged cells exactly)
# and drop the copy right after it - "Nädal 2" is the last tab in the
# workbook at this point, so that also puts the copy at the very end.
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Лист1"

# New header date.
$newSheet.Range("G4").Value2 = 40953

# Row 7: new log entry (date, start/stop time, duration, activity text).
$newSheet.Range("B7").Value2 = 43874
$newSheet.Range("C7").Value2 = 0.44444444444444442
$newSheet.Range("D7").Value2 = 0.4548611111111111
$newSheet.Range("F7").Value2 = 15
$newSheet.Range("G7").Value2 = "Kodutöö esitamine"

# Rows 8-16 had sample entries copied over from "Nädal 2" - wipe the values
# back out (formatting/styles stay untouched) so only row 7 has data.
$newSheet.Range("B8:J16").ClearContents()

# The copied-over comment on H14 no longer has any backing data - remove it
# (this also drops the sheet's legacyDrawing/VML reference).
$commentCell = $newSheet.Range("H14")
if ($commentCell.Comment -ne $null) {
    $commentCell.Comment.Delete()
}

# Make the new sheet the active tab/selection.
$newSheet.Select()
[void]$newSheet.Range("G10").Select()

# "Nädal 2" no longer holds the "last active" flag; its selection becomes a
# full-sheet (select-all) selection.
[void]$template.Range("A1:XFD1048576").Select()

$newSheet.Activate()
